# Apply the "add 2022-Q3 data" edit:
#  1. Insert a new sheet "2022-Q3" right after "总计", shifting the other
#     quarter sheets down one tab position (their content is untouched).
#  2. Insert a new summary row on "总计" for "2022-Q3" (5 funds held,
#     0.14 亿元 held value), pushing the existing rows down by one.
#  3. Populate the new "2022-Q3" sheet with its fund holdings table.

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Step 1: shift the existing "总计" quarter data (columns B-D only) down
# one row, from rows 2-8 to rows 3-9, then write the brand-new "2022-Q3"
# values into row 2. Column A is a plain 0-based row index (0,1,2,...)
# that is NOT tied to any quarter, so it is left alone for rows 2-8 and
# simply extended with the next index (7) for the newly appended row 9.
# Walk bottom-up so a row is never clobbered before it has been read.
# ---------------------------------------------------------------------
for ($r = 8; $r -ge 2; $r--) {
    $destRow = $r + 1
    $total.Range("B$r").Copy($total.Range("B$destRow"))
    $total.Range("C$r").Copy($total.Range("C$destRow"))
    $total.Range("D$r").Copy($total.Range("D$destRow"))
}

$total.Range("A8").Copy($total.Range("A9"))
$total.Range("A9").Value = 7

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 0.14

# ---------------------------------------------------------------------
# Step 2: add the new "2022-Q3" worksheet right after "总计". The other
# quarter sheets (e.g. the "2022-Q2" tab, currently 2nd) share an
# identical header-row layout/style, so that tab is used purely as a
# style template for the new header cells.
# ---------------------------------------------------------------------
$styleSrc = $wb.Worksheets.Item(2)

$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$cols = @("B", "C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $cols[$i]
    $styleSrc.Range("${col}1").Copy($q3.Range("${col}1"))
    $q3.Range("${col}1").Value = $headers[$i]
}

$rows = @(
    @(0, "000646", "华润元大量化优选混合A", "1.47", "73.62", "4.83", "0.0710", 6),
    @(1, "004260", "德邦稳盈增长灵活配置混合", "1.37", "88.81", "4.01", "0.0549", 9),
    @(2, "007827", "华润元大量化优选混合C", "0.19", "73.62", "4.83", "0.0092", 6),
    @(3, "009649", "嘉实精选平衡混合A", "0.07", "58.93", "2.28", "0.0016", 9),
    @(4, "009650", "嘉实精选平衡混合C", "0.05", "58.93", "2.28", "0.0011", 9)
)

# Columns B, D, E, F, G hold number-shaped *text* in the source data
# (fund code "000646" keeps its leading zero, percentages keep trailing
# zeros like "0.0710") - force the Text format before assigning so the
# engine does not silently coerce them to numeric cells.
$q3.Range("B2:B6").NumberFormat = "@"
$q3.Range("D2:G6").NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    $styleSrc.Range("A2").Copy($q3.Range("A$r"))
    $q3.Range("A$r").Value = $row[0]
    $q3.Range("B$r").Value = $row[1]
    $q3.Range("C$r").Value = $row[2]
    $q3.Range("D$r").Value = $row[3]
    $q3.Range("E$r").Value = $row[4]
    $q3.Range("F$r").Value = $row[5]
    $q3.Range("G$r").Value = $row[6]
    $q3.Range("H$r").Value = $row[7]
    $r++
}
